$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'270.62"
$ws.Range("D4").Value = "'6.359"
$ws.Range("D6").Value = "'3.649"
$ws.Range("D7").Value = "'6.710"
$ws.Range("D8").Value = "'1.373"
$ws.Range("D9").Value = "'0.8355"
$ws.Range("D10").Value = "'0.01379"
$ws.Range("D11").Value = "'0.1628"
$ws.Range("D12").Value = "'0.08379"
$ws.Range("D13").Value = "'0.03481"
$ws.Range("D14").Value = "'0.03140"
$ws.Range("D15").Value = "'0.09328"
$ws.Range("D16").Value = "'3.878"
$ws.Range("D17").Value = "'0.001713"
$ws.Range("D18").Value = "'0.04830"
$ws.Range("D19").Value = "'0.006178"
$ws.Range("D20").Value = "'0.001085"
$ws.Range("D21").Value = "'0.003609"
$ws.Range("D23").Value = "'3.743"
$ws.Range("D24").Value = "'2.324"
$ws.Range("D26").Value = "'0.1261"
$ws.Range("D40").Value = "'0.04682"
$ws.Range("D41").Value = "'0.006912"
$ws.Range("D42").Value = "'0.1174"
$ws.Range("D43").Value = "'0.003448"
$ws.Range("D44").Value = "'0.01139"
$ws.Range("D45").Value = "'0.00006265"
$ws.Range("D46").Value = "'0.00000000748"
$ws.Range("D47").Value = "'0.7968"
$ws.Range("D48").Value = "'0.08825"
$ws.Range("D49").Value = "'0.00001397"
$ws.Range("D50").Value = "'0.01237"
